$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the "Product" column (column C). Everything to its right shifts
#    one column to the left (PackageType moves from E to D, City from Q to P,
#    etc.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).Delete()

# ---------------------------------------------------------------------------
# 2. Append the three new header cells after the last existing column
#    (now column S = "latitude"): T=allowOpenPackages, U=feesOnConsignee,
#    V=sameDayDelivery. They are written in the order T, V, U so that the
#    shared-string table ends up in the same order as the target workbook.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,20).Value = "allowOpenPackages"
$ws.Columns.Item(20).AutoFit()

$ws.Cells.Item(1,22).Value = "sameDayDelivery"
$ws.Columns.Item(22).AutoFit()

$ws.Cells.Item(1,21).Value = "feesOnConsignee"
$ws.Columns.Item(21).AutoFit()

# ---------------------------------------------------------------------------
# 3. New data validation (Yes/No list) for the three new columns, rows 2+
# ---------------------------------------------------------------------------
$ws.Range("T2:V1048576").Validation.Add(3, 1, 1, """Yes,No""")

# ---------------------------------------------------------------------------
# 4. Cell comments explaining each new column
# ---------------------------------------------------------------------------
$cT = $ws.Range("T1").AddComment()
$cT.Text("Adham Ahmed:" + [char]10 + "+5 EGP will be added on the fees, if you don't choose Yes or No, the services you chose while registration will be applied")

$cU = $ws.Range("U1").AddComment()
$cU.Text("Adham Ahmed:" + [char]10 + "+10 EGP will be added on the fees, if you don't choose Yes or No, the services you chose while registration will be applied")

$cV = $ws.Range("V1").AddComment()
$cV.Text("Adham Ahmed:" + [char]10 + "+15 EGP will be added on the fees, if you don't choose Yes or No, the services you chose while registration will be applied")

# ---------------------------------------------------------------------------
# 5. Leave the same cell selected as in the saved workbook
# ---------------------------------------------------------------------------
$ws.Range("G8").Select()

Write-Output "Edit complete"
